# Scheduled refresh of market-price-derived columns (currentAveragePrice*,
# LevePrice*, LeveProfit*) on each crafting-job "Profits" sheet, per the
# latest market board data pull. Columns: H=currentAveragePrice,
# I=currentAveragePriceNQ, J=currentAveragePriceHQ, K=LevePriceNQ,
# L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ (col 8..14).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 934.6667
$ws.Cells.Item(33, 9).Value = 989
$ws.Cells.Item(33, 11).Value = 989
$ws.Cells.Item(33, 13).Value = -760
$ws.Cells.Item(58, 8).Value = 2228.1428
$ws.Cells.Item(58, 9).Value = 1649.25
$ws.Cells.Item(58, 11).Value = 4947.75
$ws.Cells.Item(58, 13).Value = -4797.75
$ws.Cells.Item(137, 8).Value = 2138116
$ws.Cells.Item(137, 9).Value = 4386875.5
$ws.Cells.Item(137, 10).Value = 1794.2
$ws.Cells.Item(137, 11).Value = 13160626.5
$ws.Cells.Item(137, 12).Value = 5382.6
$ws.Cells.Item(137, 13).Value = -13158076.5
$ws.Cells.Item(137, 14).Value = -10482.6
$ws.Cells.Item(138, 8).Value = 3783.7666
$ws.Cells.Item(138, 9).Value = 6824.857
$ws.Cells.Item(138, 10).Value = 3527.289
$ws.Cells.Item(138, 11).Value = 20474.571
$ws.Cells.Item(138, 12).Value = 10581.867
$ws.Cells.Item(138, 13).Value = -15334.571
$ws.Cells.Item(138, 14).Value = -20861.867
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(23, 8).Value = 8108.826
$ws.Cells.Item(23, 10).Value = 4840.773
$ws.Cells.Item(23, 12).Value = 4840.773
$ws.Cells.Item(23, 14).Value = -5358.773
$ws.Cells.Item(44, 8).Value = 27399.2
$ws.Cells.Item(44, 9).Value = 17000
$ws.Cells.Item(44, 11).Value = 17000
$ws.Cells.Item(44, 13).Value = -16512
$ws.Cells.Item(55, 8).Value = 27499.166
$ws.Cells.Item(55, 9).Value = 15000
$ws.Cells.Item(55, 11).Value = 15000
$ws.Cells.Item(55, 13).Value = -14685
$ws.Cells.Item(61, 8).Value = 11113680
$ws.Cells.Item(61, 9).Value = 13334976
$ws.Cells.Item(61, 10).Value = 7200
$ws.Cells.Item(61, 11).Value = 13334976
$ws.Cells.Item(61, 12).Value = 7200
$ws.Cells.Item(61, 13).Value = -13334764
$ws.Cells.Item(61, 14).Value = -7624
$ws.Cells.Item(104, 8).Value = 75890
$ws.Cells.Item(104, 10).Value = 75890
$ws.Cells.Item(104, 12).Value = 75890
$ws.Cells.Item(104, 14).Value = -82878
$ws.Cells.Item(136, 8).Value = 11113680
$ws.Cells.Item(136, 9).Value = 13334976
$ws.Cells.Item(136, 10).Value = 7200
$ws.Cells.Item(136, 11).Value = 40004928
$ws.Cells.Item(136, 12).Value = 21600
$ws.Cells.Item(136, 13).Value = -40002378
$ws.Cells.Item(136, 14).Value = -26700
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 24391702
$ws.Cells.Item(20, 9).Value = 1364.92
$ws.Cells.Item(20, 10).Value = 62501600
$ws.Cells.Item(20, 11).Value = 1364.92
$ws.Cells.Item(20, 12).Value = 62501600
$ws.Cells.Item(20, 13).Value = -1117.92
$ws.Cells.Item(20, 14).Value = -62502094
$ws.Cells.Item(82, 8).Value = 15564.353
$ws.Cells.Item(82, 9).Value = 3475.125
$ws.Cells.Item(82, 10).Value = 26310.334
$ws.Cells.Item(82, 11).Value = 3475.125
$ws.Cells.Item(82, 12).Value = 26310.334
$ws.Cells.Item(82, 13).Value = -3092.125
$ws.Cells.Item(82, 14).Value = -27076.334
$ws.Cells.Item(85, 8).Value = 15564.353
$ws.Cells.Item(85, 9).Value = 3475.125
$ws.Cells.Item(85, 10).Value = 26310.334
$ws.Cells.Item(85, 11).Value = 3475.125
$ws.Cells.Item(85, 12).Value = 26310.334
$ws.Cells.Item(85, 13).Value = -2149.125
$ws.Cells.Item(85, 14).Value = -28962.334
$ws.Cells.Item(116, 8).Value = 25870.75
$ws.Cells.Item(116, 10).Value = 25870.75
$ws.Cells.Item(116, 12).Value = 25870.75
$ws.Cells.Item(116, 14).Value = -35048.75
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(41, 8).Value = 13995
$ws.Cells.Item(41, 10).Value = 15151.429
$ws.Cells.Item(41, 12).Value = 15151.429
$ws.Cells.Item(41, 14).Value = -16007.429
$ws.Cells.Item(50, 8).Value = 11999.375
$ws.Cells.Item(50, 10).Value = 12285
$ws.Cells.Item(50, 12).Value = 12285
$ws.Cells.Item(50, 14).Value = -13535
$ws.Cells.Item(51, 8).Value = 16666
$ws.Cells.Item(51, 10).Value = 17999.2
$ws.Cells.Item(51, 12).Value = 17999.2
$ws.Cells.Item(51, 14).Value = -19471.2
$ws.Cells.Item(58, 8).Value = 1186.55
$ws.Cells.Item(58, 9).Value = 884.8461
$ws.Cells.Item(58, 10).Value = 1746.8572
$ws.Cells.Item(58, 11).Value = 884.8461
$ws.Cells.Item(58, 12).Value = 1746.8572
$ws.Cells.Item(58, 13).Value = -681.8461
$ws.Cells.Item(58, 14).Value = -2152.8572
$ws.Cells.Item(59, 8).Value = 19999.166
$ws.Cells.Item(59, 10).Value = 19999.166
$ws.Cells.Item(59, 12).Value = 19999.166
$ws.Cells.Item(59, 14).Value = -22289.166
$ws.Cells.Item(60, 8).Value = 10975.765
$ws.Cells.Item(60, 9).Value = 1036
$ws.Cells.Item(60, 10).Value = 15117.333
$ws.Cells.Item(60, 11).Value = 1036
$ws.Cells.Item(60, 12).Value = 15117.333
$ws.Cells.Item(60, 13).Value = -525
$ws.Cells.Item(60, 14).Value = -16139.333
$ws.Cells.Item(61, 8).Value = 16666
$ws.Cells.Item(61, 10).Value = 17999.2
$ws.Cells.Item(61, 12).Value = 17999.2
$ws.Cells.Item(61, 14).Value = -18695.2
$ws.Cells.Item(68, 8).Value = 22999.46
$ws.Cells.Item(68, 10).Value = 22999.46
$ws.Cells.Item(68, 12).Value = 22999.46
$ws.Cells.Item(68, 14).Value = -24497.46
$ws.Cells.Item(71, 8).Value = 22999.46
$ws.Cells.Item(71, 10).Value = 22999.46
$ws.Cells.Item(71, 12).Value = 68998.38
$ws.Cells.Item(71, 14).Value = -76486.38
$ws.Cells.Item(74, 8).Value = 19599.4
$ws.Cells.Item(74, 10).Value = 19599.4
$ws.Cells.Item(74, 12).Value = 19599.4
$ws.Cells.Item(74, 14).Value = -21347.4
$ws.Cells.Item(77, 8).Value = 19599.4
$ws.Cells.Item(77, 10).Value = 19599.4
$ws.Cells.Item(77, 12).Value = 58798.2
$ws.Cells.Item(77, 14).Value = -67534.20000000001
$ws.Cells.Item(132, 8).Value = 30305772
$ws.Cells.Item(132, 9).Value = 45457210
$ws.Cells.Item(132, 10).Value = 15154335
$ws.Cells.Item(132, 11).Value = 136371630
$ws.Cells.Item(132, 12).Value = 45463005
$ws.Cells.Item(132, 13).Value = -136369100
$ws.Cells.Item(132, 14).Value = -45468065
$ws.Cells.Item(136, 8).Value = 1186.55
$ws.Cells.Item(136, 9).Value = 884.8461
$ws.Cells.Item(136, 10).Value = 1746.8572
$ws.Cells.Item(136, 11).Value = 2654.5383
$ws.Cells.Item(136, 12).Value = 5240.571599999999
$ws.Cells.Item(136, 13).Value = -104.5383000000002
$ws.Cells.Item(136, 14).Value = -10340.5716
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 924.7442
$ws.Cells.Item(5, 9).Value = 717
$ws.Cells.Item(5, 10).Value = 1163.65
$ws.Cells.Item(5, 11).Value = 2151
$ws.Cells.Item(5, 12).Value = 3490.95
$ws.Cells.Item(5, 13).Value = -2039
$ws.Cells.Item(5, 14).Value = -3714.95
$ws.Cells.Item(135, 8).Value = 924.7442
$ws.Cells.Item(135, 9).Value = 717
$ws.Cells.Item(135, 10).Value = 1163.65
$ws.Cells.Item(135, 11).Value = 6453
$ws.Cells.Item(135, 12).Value = 10472.85
$ws.Cells.Item(135, 13).Value = -3918
$ws.Cells.Item(135, 14).Value = -15542.85
$ws.Cells.Item(140, 8).Value = 1558.7106
$ws.Cells.Item(140, 9).Value = 1011.11536
$ws.Cells.Item(140, 10).Value = 2745.1667
$ws.Cells.Item(140, 11).Value = 3033.34608
$ws.Cells.Item(140, 12).Value = 8235.500100000001
$ws.Cells.Item(140, 13).Value = 2146.65392
$ws.Cells.Item(140, 14).Value = -18595.5001
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 4015.25
$ws.Cells.Item(7, 9).Value = 4551.3076
$ws.Cells.Item(7, 10).Value = 3381.7273
$ws.Cells.Item(7, 11).Value = 4551.3076
$ws.Cells.Item(7, 12).Value = 3381.7273
$ws.Cells.Item(7, 13).Value = -4439.3076
$ws.Cells.Item(7, 14).Value = -3605.7273
$ws.Cells.Item(46, 8).Value = 599.875
$ws.Cells.Item(46, 9).Value = 649.75
$ws.Cells.Item(46, 10).Value = 550
$ws.Cells.Item(46, 11).Value = 649.75
$ws.Cells.Item(46, 12).Value = 550
$ws.Cells.Item(46, 13).Value = -461.75
$ws.Cells.Item(46, 14).Value = -926
$ws.Cells.Item(55, 8).Value = 601.8333
$ws.Cells.Item(55, 9).Value = 243.66667
$ws.Cells.Item(55, 11).Value = 243.66667
$ws.Cells.Item(55, 13).Value = -70.66667000000001
$ws.Cells.Item(93, 8).Value = 9556.691999999999
$ws.Cells.Item(93, 9).Value = 17985
$ws.Cells.Item(93, 10).Value = 2332.4285
$ws.Cells.Item(93, 11).Value = 17985
$ws.Cells.Item(93, 12).Value = 2332.4285
$ws.Cells.Item(93, 13).Value = -16737
$ws.Cells.Item(93, 14).Value = -4828.4285
$ws.Cells.Item(112, 8).Value = 84188.5
$ws.Cells.Item(112, 10).Value = 84188.5
$ws.Cells.Item(112, 12).Value = 84188.5
$ws.Cells.Item(112, 14).Value = -87142.5
$ws.Cells.Item(122, 8).Value = 3946.8157
$ws.Cells.Item(122, 9).Value = 2576.389
$ws.Cells.Item(122, 10).Value = 5180.2
$ws.Cells.Item(122, 11).Value = 7729.167
$ws.Cells.Item(122, 12).Value = 15540.6
$ws.Cells.Item(122, 13).Value = -5279.167
$ws.Cells.Item(122, 14).Value = -20440.6
$ws.Cells.Item(126, 8).Value = 4015.25
$ws.Cells.Item(126, 9).Value = 4551.3076
$ws.Cells.Item(126, 10).Value = 3381.7273
$ws.Cells.Item(126, 11).Value = 13653.9228
$ws.Cells.Item(126, 12).Value = 10145.1819
$ws.Cells.Item(126, 13).Value = -11183.9228
$ws.Cells.Item(126, 14).Value = -15085.1819
$ws.Cells.Item(137, 8).Value = 39950
$ws.Cells.Item(137, 10).Value = 39950
$ws.Cells.Item(137, 12).Value = 39950
$ws.Cells.Item(137, 14).Value = -50150
$ws.Cells.Item(140, 8).Value = 61125
$ws.Cells.Item(140, 10).Value = 61125
$ws.Cells.Item(140, 12).Value = 61125
$ws.Cells.Item(140, 14).Value = -71485
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(110, 8).Value = 98644
$ws.Cells.Item(110, 10).Value = 98644
$ws.Cells.Item(110, 12).Value = 98644
$ws.Cells.Item(110, 14).Value = -106824
